$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes the existing rows 3-15 down to 4-16),
# mirroring a new weekly price record being added at the top of the data.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's record.
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Femacal de La Calera"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44959
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101004
$ws.Range("J3").Value = "Frambuesa"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 40
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("Q3").Value = "$/bandeja 2 kilos"
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 3500
$ws.Range("T3").Value = 2
